$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of trade data (row 5)
$ws.Range("A5").Value = 42636.593078703707
$ws.Range("B5").Value = $true
$ws.Range("C5").Value = 9979.36
$ws.Range("D5").Value = 9951.5
$ws.Range("E5").Value = 313.07
$ws.Range("F5").Value = 314.81
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 0.56000000000000005
$ws.Range("I5").Value = $false

# Match the date-formatted style used by the rows above (A/G columns)
$ws.Range("A3").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("G3").Copy()
$ws.Range("G5").PasteSpecial(-4122)
